$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), copying the header style from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2..51 for columns I (I0) and J (IF)
$data = @(
    @(2, 9, 9),
    @(3, 6, 6),
    @(4, 7, 7),
    @(5, 7, 7),
    @(6, 6, 6),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 8, 8),
    @(11, 7, 7),
    @(12, 7, 7),
    @(13, 8, 8),
    @(14, 4, 4),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 8, 8),
    @(18, 6, 7),
    @(19, 8, 8),
    @(20, 7, 8),
    @(21, 8, 8),
    @(22, 5, 7),
    @(23, 8, 8),
    @(24, 4, 5),
    @(25, 7, 7),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 8, 9),
    @(30, 7, 8),
    @(31, 7, 7),
    @(32, 6, 7),
    @(33, 7, 7),
    @(34, 8, 8),
    @(35, 5, 6),
    @(36, 9, 9),
    @(37, 5, 6),
    @(38, 6, 7),
    @(39, 7, 8),
    @(40, 6, 8),
    @(41, 6, 7),
    @(42, 7, 8),
    @(43, 6, 7),
    @(44, 7, 9),
    @(45, 10, 10),
    @(46, 7, 7),
    @(47, 4, 5),
    @(48, 4, 4),
    @(49, 4, 4),
    @(50, 5, 5),
    @(51, 6, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
